# Insert a new data row at row 310 (pushing the existing rows 310-388 down
# to 311-389) and populate it with the new Cilantro price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 310..388 down by inserting a new row at position 310.
$ws.Rows.Item(310).Insert()

# Populate the newly inserted row 310 with the new record's values.
$ws.Cells.Item(310, 1).Value = 4
$ws.Cells.Item(310, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(310, 3).Value = "Los Lagos"
$ws.Cells.Item(310, 4).Value = 44943
$ws.Cells.Item(310, 5).Value = 10
$ws.Cells.Item(310, 6).Value = 100112040
$ws.Cells.Item(310, 7).Value = "Cilantro"
$ws.Cells.Item(310, 8).Value = "Sin especificar"
$ws.Cells.Item(310, 9).Value = "Primera"
$ws.Cells.Item(310, 10).Value = 160
$ws.Cells.Item(310, 11).Value = 8000
$ws.Cells.Item(310, 12).Value = 9000
$ws.Cells.Item(310, 13).Value = 8500
$ws.Cells.Item(310, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(310, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(310, 16).Value = 4250
$ws.Cells.Item(310, 17).Value = 2
$ws.Cells.Item(310, 18).Value = "Hortaliza"
